$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  @(2287.81,2731.01,4300.83,4977.32,4791.22,4332.28,3803.81,3381.93,3137.76,3023.37,2975.24,2955.91,2948.3,2945.33),
  @(2532.19,2789.23,4172.95,4809.86,4552.18,4048.08,3572.22,3265.87,3113.5,3047.56,3020.78,3010.18,3006.03,3004.41),
  @(2893.4,2987.66,4142.43,4700.02,4367.59,3881.49,3515.53,3317.8,3228.8,3192.02,3177.38,3171.62,3169.38,3168.5),
  @(3369.4,3339.81,4237.68,4688.55,4317.99,3917.51,3672.24,3555.28,3505.74,3485.8,3477.94,3474.87,3473.67,3473.2),
  @(3941.11,3821.79,4467.43,4828.62,4482.67,4195.31,4046.04,3980.57,3953.84,3943.25,3939.1,3937.48,3936.85,3936.6),
  @(4554.83,4357.36,4803.98,5154.9,4880.32,4685.38,4595,4557.48,4542.51,4536.64,4534.34,4533.45,4533.1,4532.96),
  @(5111.12,4831.74,5180.17,5636.09,5443.29,5294.61,5226.58,5198.58,5187.47,5183.11,5181.41,5180.74,5180.49,5180.39),
  @(5496.16,5147.12,5511.97,6171.5,6056.89,5907.82,5828.78,5794.08,5779.92,5774.3,5772.1,5771.25,5770.91,5770.78),
  @(5646.88,5276.9,5734.17,6647.43,6624.01,6445.17,6324.83,6266.41,6241.56,6231.55,6227.6,6226.06,6225.46,6225.22),
  @(5586.86,5258.22,5821.5,6990.73,7093.56,6885.15,6699.09,6597.39,6551.91,6533.2,6525.76,6522.84,6521.71,6521.26),
  @(5397.08,5144.03,5784.35,7180.28,7447.84,7242.61,6980.58,6814.18,6734.46,6700.7,6687.13,6681.77,6679.68,6678.87),
  @(5156.67,4974.1,5651.78,7231.12,7685.52,7535.7,7212.02,6962.67,6830.99,6772.73,6748.86,6739.39,6735.67,6734.23),
  @(4910.58,4773.46,5456.34,7174.01,7816.02,7769.7,7425.71,7088.39,6884.69,6788.39,6747.78,6731.47,6725.05,6722.54),
  @(4674.68,4558.38,5226.06,7041.11,7857.18,7941.18,7631.76,7225.4,6933.3,6781.22,6714.13,6686.67,6675.77,6671.5),
  @(4451.65,4340.2,4981.97,6859.15,7830.44,8049.27,7820.62,7388.93,7006.46,6778.93,6671.45,6626.1,6607.87,6600.69),
  @(4240.77,4126.83,4738.39,6647.98,7755.78,8100.52,7976.44,7572.67,7122.23,6804.45,6638.59,6565.2,6535.1,6523.15)
)

for ($ri = 0; $ri -lt $data.Count; $ri++) {
  $rowVals = $data[$ri]
  for ($ci = 0; $ci -lt $rowVals.Count; $ci++) {
    $ws.Cells.Item(4 + $ri, 2 + $ci).Value = $rowVals[$ci]
  }
}

$wb.Save()